$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = -0.08634205848777793
$ws.Range("J2").Value = 0.1244354596818399
$ws.Range("K2").Value = -0.5778595736595014
$ws.Range("L2").Value = 2.300016341032164

# Row 12
$ws.Range("I12").Value = 0.1451106724099447
$ws.Range("J12").Value = 0.03860845840830766
$ws.Range("K12").Value = -0.719497518369239
$ws.Range("L12").Value = 2.141425594417664

# Row 13
$ws.Range("I13").Value = 0.03812022684306444
$ws.Range("J13").Value = 0.08940849255325421
$ws.Range("K13").Value = -0.7000752725432544
$ws.Range("L13").Value = 2.086799292840594

# Row 14
$ws.Range("I14").Value = 0.003963586090688369
$ws.Range("J14").Value = 0.08363590740485342
$ws.Range("K14").Value = -0.3880079538026978
$ws.Range("L14").Value = 1.843662676263534

# Row 16
$ws.Range("I16").Value = -0.000265022413729249
$ws.Range("J16").Value = 0.0733065868884322
$ws.Range("K16").Value = -0.3671219725154395
$ws.Range("L16").Value = 2.075485821741652

# Row 18
$ws.Range("I18").Value = -0.2782813428863373
$ws.Range("J18").Value = 0.1409092920061238
$ws.Range("K18").Value = 0.1129486708220114
$ws.Range("L18").Value = 1.875111306265327
